# Generate Report for Handoff
# Adds a new row (the 9b510dd0-... handback file) to the Overview, zh-cn
# and de-de tables/sheets.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d000a92316455503912bb528b49065fb34f9d7f/e2e/9b510dd0-dead-46e1-873c-08c665664ef7.md"

# ---------------------------------------------------------------------
# Overview sheet: File Name | Path And Name | Extension | Publish URL |
#                 zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "9b510dd0-dead-46e1-873c-08c665664ef7.md"
$wsOverview.Range("B4").Value = "e2e\9b510dd0-dead-46e1-873c-08c665664ef7.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = "'"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-10-24 10:05:47"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $baseUrl, "", "", "e2e\9b510dd0-dead-46e1-873c-08c665664ef7.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet (16 columns)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = "9b510dd0-dead-46e1-873c-08c665664ef7.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'False"
$wsZhCn.Range("G4").Value = "9b510dd0-dead-46e1-873c-08c665664ef7.be478683d9cdcad813448382cb1f2f0bfc81d2f1.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-10-24 10:05:35"
$wsZhCn.Range("I4").Value = "'"
$wsZhCn.Range("J4").Value = "'"
$wsZhCn.Range("K4").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $baseUrl, "", "", "9b510dd0-dead-46e1-873c-08c665664ef7.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet (16 columns)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = "9b510dd0-dead-46e1-873c-08c665664ef7.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'False"
$wsDeDe.Range("G4").Value = "9b510dd0-dead-46e1-873c-08c665664ef7.be478683d9cdcad813448382cb1f2f0bfc81d2f1.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-10-24 10:05:47"
$wsDeDe.Range("I4").Value = "'"
$wsDeDe.Range("J4").Value = "'"
$wsDeDe.Range("K4").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $baseUrl, "", "", "9b510dd0-dead-46e1-873c-08c665664ef7.md") | Out-Null

Write-Host "Added handback row for 9b510dd0-dead-46e1-873c-08c665664ef7.md to Overview, zh-cn and de-de sheets."
